# Added another screen for quick interactivity testing
# Adds a new "Articles" sub-section (column D) to Sheet1 describing an
# additional screen concept, and moves the active-cell selection to D32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D13").Value = "Screen Transitions (ON/OFF)"

$ws.Range("D20").Value = "Articles"
$ws.Range("D20").Font.Bold = $true

$ws.Range("D21").Value = "Flash is Dead, Long Live HTML5!"

$ws.Range("D24").Value = "Hybrid Adventures:"
$ws.Range("D25").Value = "SPA"
$ws.Range("D26").Value = "Native vs. Mobile"
$ws.Range("D27").Value = "Hybrid Trend"

$ws.Range("D22").Value = "iOS WebView, WebGL defaultly on, …"

$ws.Range("D28").Value = "Canvas vs. HTML/CSS layout inconsistencies"
$ws.Range("D29").Value = "Canvas vs. HTML/CSS performance"

$ws.Range("D32").Select()
